$d = $word.ActiveDocument

# Date line update
$d.Content.Find.Execute("2024-08-25 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-26 Monday", 2)

# Table of multiplication problems. The value "98x96=9408" appears twice with
# different replacements, so address those two cells directly via the table
# before doing the remaining global (unique) replacements.
$t = $d.Tables.Item(1)
$t.Cell(1, 4).Range.Text = "23x70=1610".Replace("x", [char]0x00D7)
$t.Cell(10, 5).Range.Text = "16x35=560".Replace("x", [char]0x00D7)

# Remaining unique replacements (safe to do as document-wide Find/Replace)
$pairs = @(
    @("97x58=5626", "49x66=3234"),
    @("33x22=726", "36x63=2268"),
    @("88x53=4664", "33x45=1485"),
    @("13x89=1157", "89x67=5963"),
    @("61x13=793", "84x66=5544"),
    @("14x86=1204", "31x18=558"),
    @("82x91=7462", "25x26=650"),
    @("69x34=2346", "31x62=1922"),
    @("19x86=1634", "50x22=1100"),
    @("86x92=7912", "82x96=7872"),
    @("81x73=5913", "94x88=8272"),
    @("32x59=1888", "40x24=960"),
    @("73x56=4088", "97x47=4559"),
    @("38x92=3496", "51x30=1530"),
    @("32x69=2208", "46x23=1058"),
    @("26x32=832", "67x86=5762"),
    @("87x78=6786", "44x80=3520"),
    @("36x26=936", "49x69=3381"),
    @("16x70=1120", "74x58=4292"),
    @("52x39=2028", "86x44=3784"),
    @("48x22=1056", "55x57=3135"),
    @("85x97=8245", "97x79=7663"),
    @("21x24=504", "57x51=2907")
)

foreach ($pair in $pairs) {
    $old = $pair[0].Replace("x", [char]0x00D7)
    $new = $pair[1].Replace("x", [char]0x00D7)
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
